$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-24 Wednesday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-12-25 Thursday", 2) | Out-Null
$d.Content.Find.Execute("2+14=16", $true, $true, $false, $false, $false, $true, 1, $false, "27+24=51", 2) | Out-Null
$d.Content.Find.Execute("68-44=24", $true, $true, $false, $false, $false, $true, 1, $false, "12+67=79", 2) | Out-Null
$d.Content.Find.Execute("61+28=89", $true, $true, $false, $false, $false, $true, 1, $false, "70-30=40", 2) | Out-Null
$d.Content.Find.Execute("22+20=42", $true, $true, $false, $false, $false, $true, 1, $false, "38+5=43", 2) | Out-Null
$d.Content.Find.Execute("60+26=86", $true, $true, $false, $false, $false, $true, 1, $false, "31+40=71", 2) | Out-Null
$d.Content.Find.Execute("65+6=71", $true, $true, $false, $false, $false, $true, 1, $false, "72-33=39", 2) | Out-Null
$d.Content.Find.Execute("15+23=38", $true, $true, $false, $false, $false, $true, 1, $false, "1+48=49", 2) | Out-Null
$d.Content.Find.Execute("9-3=6", $true, $true, $false, $false, $false, $true, 1, $false, "76-8=68", 2) | Out-Null
$d.Content.Find.Execute("31+12=43", $true, $true, $false, $false, $false, $true, 1, $false, "37+8=45", 2) | Out-Null
$d.Content.Find.Execute("18+25=43", $true, $true, $false, $false, $false, $true, 1, $false, "93-42=51", 2) | Out-Null
$d.Content.Find.Execute("26+34=60", $true, $true, $false, $false, $false, $true, 1, $false, "78-19=59", 2) | Out-Null
$d.Content.Find.Execute("94-38=56", $true, $true, $false, $false, $false, $true, 1, $false, "73+7=80", 2) | Out-Null
$d.Content.Find.Execute("34+21=55", $true, $true, $false, $false, $false, $true, 1, $false, "86-11=75", 2) | Out-Null
$d.Content.Find.Execute("8+56=64", $true, $true, $false, $false, $false, $true, 1, $false, "61-28=33", 2) | Out-Null
$d.Content.Find.Execute("76-31=45", $true, $true, $false, $false, $false, $true, 1, $false, "1+62=63", 2) | Out-Null
$d.Content.Find.Execute("98-30=68", $true, $true, $false, $false, $false, $true, 1, $false, "14+79=93", 2) | Out-Null
$d.Content.Find.Execute("65-48=17", $true, $true, $false, $false, $false, $true, 1, $false, "39+50=89", 2) | Out-Null
$d.Content.Find.Execute("18+67=85", $true, $true, $false, $false, $false, $true, 1, $false, "11-5=6", 2) | Out-Null
$d.Content.Find.Execute("92-51=41", $true, $true, $false, $false, $false, $true, 1, $false, "54-23=31", 2) | Out-Null
$d.Content.Find.Execute("70-3=67", $true, $true, $false, $false, $false, $true, 1, $false, "40+6=46", 2) | Out-Null
$d.Content.Find.Execute("44+20=64", $true, $true, $false, $false, $false, $true, 1, $false, "74+18=92", 2) | Out-Null
$d.Content.Find.Execute("13-2=11", $true, $true, $false, $false, $false, $true, 1, $false, "13-13=0", 2) | Out-Null
$d.Content.Find.Execute("66-28=38", $true, $true, $false, $false, $false, $true, 1, $false, "77-22=55", 2) | Out-Null
$d.Content.Find.Execute("37-0=37", $true, $true, $false, $false, $false, $true, 1, $false, "68-0=68", 2) | Out-Null
$d.Content.Find.Execute("53-4=49", $true, $true, $false, $false, $false, $true, 1, $false, "81-80=1", 2) | Out-Null
$d.Content.Find.Execute("39+25=64", $true, $true, $false, $false, $false, $true, 1, $false, "44-32=12", 2) | Out-Null
$d.Content.Find.Execute("14+22=36", $true, $true, $false, $false, $false, $true, 1, $false, "57+16=73", 2) | Out-Null
$d.Content.Find.Execute("70-43=27", $true, $true, $false, $false, $false, $true, 1, $false, "93-60=33", 2) | Out-Null
$d.Content.Find.Execute("97-24=73", $true, $true, $false, $false, $false, $true, 1, $false, "93-84=9", 2) | Out-Null
$d.Content.Find.Execute("18+55=73", $true, $true, $false, $false, $false, $true, 1, $false, "9+66=75", 2) | Out-Null
$d.Content.Find.Execute("79-44=35", $true, $true, $false, $false, $false, $true, 1, $false, "94-91=3", 2) | Out-Null
$d.Content.Find.Execute("66-25=41", $true, $true, $false, $false, $false, $true, 1, $false, "37+13=50", 2) | Out-Null
$d.Content.Find.Execute("6+81=87", $true, $true, $false, $false, $false, $true, 1, $false, "84-41=43", 2) | Out-Null
$d.Content.Find.Execute("70-40=30", $true, $true, $false, $false, $false, $true, 1, $false, "72+6=78", 2) | Out-Null
$d.Content.Find.Execute("10+70=80", $true, $true, $false, $false, $false, $true, 1, $false, "55+13=68", 2) | Out-Null
$d.Content.Find.Execute("36+56=92", $true, $true, $false, $false, $false, $true, 1, $false, "0+53=53", 2) | Out-Null
$d.Content.Find.Execute("65-44=21", $true, $true, $false, $false, $false, $true, 1, $false, "88+7=95", 2) | Out-Null
$d.Content.Find.Execute("64-56=8", $true, $true, $false, $false, $false, $true, 1, $false, "72-14=58", 2) | Out-Null
$d.Content.Find.Execute("9+79=88", $true, $true, $false, $false, $false, $true, 1, $false, "19+2=21", 2) | Out-Null
$d.Content.Find.Execute("56+29=85", $true, $true, $false, $false, $false, $true, 1, $false, "81-11=70", 2) | Out-Null
$d.Content.Find.Execute("27+18=45", $true, $true, $false, $false, $false, $true, 1, $false, "73+23=96", 2) | Out-Null
$d.Content.Find.Execute("58+14=72", $true, $true, $false, $false, $false, $true, 1, $false, "9+1=10", 2) | Out-Null
$d.Content.Find.Execute("58+40=98", $true, $true, $false, $false, $false, $true, 1, $false, "11+3=14", 2) | Out-Null
$d.Content.Find.Execute("39-34=5", $true, $true, $false, $false, $false, $true, 1, $false, "5+92=97", 2) | Out-Null
$d.Content.Find.Execute("45+52=97", $true, $true, $false, $false, $false, $true, 1, $false, "93-33=60", 2) | Out-Null
$d.Content.Find.Execute("4+7=11", $true, $true, $false, $false, $false, $true, 1, $false, "85+1=86", 2) | Out-Null
$d.Content.Find.Execute("55-23=32", $true, $true, $false, $false, $false, $true, 1, $false, "1+97=98", 2) | Out-Null
$d.Content.Find.Execute("14+31=45", $true, $true, $false, $false, $false, $true, 1, $false, "95-89=6", 2) | Out-Null
$d.Content.Find.Execute("54+22=76", $true, $true, $false, $false, $false, $true, 1, $false, "41-13=28", 2) | Out-Null
$d.Content.Find.Execute("62+16=78", $true, $true, $false, $false, $false, $true, 1, $false, "58+1=59", 2) | Out-Null
$d.Content.Find.Execute("68+12=80", $true, $true, $false, $false, $false, $true, 1, $false, "93-67=26", 2) | Out-Null
$d.Content.Find.Execute("2+76=78", $true, $true, $false, $false, $false, $true, 1, $false, "82+7=89", 2) | Out-Null
$d.Content.Find.Execute("82-2=80", $true, $true, $false, $false, $false, $true, 1, $false, "69-12=57", 2) | Out-Null
$d.Content.Find.Execute("66+18=84", $true, $true, $false, $false, $false, $true, 1, $false, "50-25=25", 2) | Out-Null
$d.Content.Find.Execute("70+29=99", $true, $true, $false, $false, $false, $true, 1, $false, "63-39=24", 2) | Out-Null
$d.Content.Find.Execute("26+12=38", $true, $true, $false, $false, $false, $true, 1, $false, "66-45=21", 2) | Out-Null
$d.Content.Find.Execute("87-58=29", $true, $true, $false, $false, $false, $true, 1, $false, "62-6=56", 2) | Out-Null
$d.Content.Find.Execute("68-58=10", $true, $true, $false, $false, $false, $true, 1, $false, "48+10=58", 2) | Out-Null
$d.Content.Find.Execute("31+24=55", $true, $true, $false, $false, $false, $true, 1, $false, "36+51=87", 2) | Out-Null
$d.Content.Find.Execute("5+19=24", $true, $true, $false, $false, $false, $true, 1, $false, "18+13=31", 2) | Out-Null
$d.Content.Find.Execute("51-28=23", $true, $true, $false, $false, $false, $true, 1, $false, "58-38=20", 2) | Out-Null
$d.Content.Find.Execute("31+49=80", $true, $true, $false, $false, $false, $true, 1, $false, "84-79=5", 2) | Out-Null
$d.Content.Find.Execute("58-4=54", $true, $true, $false, $false, $false, $true, 1, $false, "13+27=40", 2) | Out-Null
$d.Content.Find.Execute("36+20=56", $true, $true, $false, $false, $false, $true, 1, $false, "50+26=76", 2) | Out-Null
$d.Content.Find.Execute("21+59=80", $true, $true, $false, $false, $false, $true, 1, $false, "98-65=33", 2) | Out-Null
$d.Content.Find.Execute("38+7=45", $true, $true, $false, $false, $false, $true, 1, $false, "53+37=90", 2) | Out-Null
$d.Content.Find.Execute("61-18=43", $true, $true, $false, $false, $false, $true, 1, $false, "94-20=74", 2) | Out-Null
$d.Content.Find.Execute("79-4=75", $true, $true, $false, $false, $false, $true, 1, $false, "60-41=19", 2) | Out-Null
$d.Content.Find.Execute("37+47=84", $true, $true, $false, $false, $false, $true, 1, $false, "33+40=73", 2) | Out-Null
$d.Content.Find.Execute("7+37=44", $true, $true, $false, $false, $false, $true, 1, $false, "83-75=8", 2) | Out-Null
$d.Content.Find.Execute("67-63=4", $true, $true, $false, $false, $false, $true, 1, $false, "45+22=67", 2) | Out-Null
$d.Content.Find.Execute("93+4=97", $true, $true, $false, $false, $false, $true, 1, $false, "82-33=49", 2) | Out-Null
$d.Content.Find.Execute("7+4=11", $true, $true, $false, $false, $false, $true, 1, $false, "52-47=5", 2) | Out-Null
$d.Content.Find.Execute("53-51=2", $true, $true, $false, $false, $false, $true, 1, $false, "58-50=8", 2) | Out-Null
$d.Content.Find.Execute("0+38=38", $true, $true, $false, $false, $false, $true, 1, $false, "47-19=28", 2) | Out-Null
$d.Content.Find.Execute("9-1=8", $true, $true, $false, $false, $false, $true, 1, $false, "63+24=87", 2) | Out-Null
$d.Content.Find.Execute("99-93=6", $true, $true, $false, $false, $false, $true, 1, $false, "4+94=98", 2) | Out-Null
$d.Content.Find.Execute("63-34=29", $true, $true, $false, $false, $false, $true, 1, $false, "21+6=27", 2) | Out-Null
$d.Content.Find.Execute("74-16=58", $true, $true, $false, $false, $false, $true, 1, $false, "88-37=51", 2) | Out-Null
$d.Content.Find.Execute("18+27=45", $true, $true, $false, $false, $false, $true, 1, $false, "12+78=90", 2) | Out-Null
$d.Content.Find.Execute("39+59=98", $true, $true, $false, $false, $false, $true, 1, $false, "58+5=63", 2) | Out-Null
$d.Content.Find.Execute("37+25=62", $true, $true, $false, $false, $false, $true, 1, $false, "85-16=69", 2) | Out-Null
$d.Content.Find.Execute("40+23=63", $true, $true, $false, $false, $false, $true, 1, $false, "70+15=85", 2) | Out-Null
$d.Content.Find.Execute("53-2=51", $true, $true, $false, $false, $false, $true, 1, $false, "13+13=26", 2) | Out-Null
$d.Content.Find.Execute("62-29=33", $true, $true, $false, $false, $false, $true, 1, $false, "65-62=3", 2) | Out-Null
$d.Content.Find.Execute("73-8=65", $true, $true, $false, $false, $false, $true, 1, $false, "86-43=43", 2) | Out-Null
$d.Content.Find.Execute("60-25=35", $true, $true, $false, $false, $false, $true, 1, $false, "33+34=67", 2) | Out-Null
$d.Content.Find.Execute("0+1=1", $true, $true, $false, $false, $false, $true, 1, $false, "4+57=61", 2) | Out-Null
$d.Content.Find.Execute("64+17=81", $true, $true, $false, $false, $false, $true, 1, $false, "67-20=47", 2) | Out-Null
$d.Content.Find.Execute("10+59=69", $true, $true, $false, $false, $false, $true, 1, $false, "34-13=21", 2) | Out-Null
$d.Content.Find.Execute("57-31=26", $true, $true, $false, $false, $false, $true, 1, $false, "78-4=74", 2) | Out-Null
$d.Content.Find.Execute("38-0=38", $true, $true, $false, $false, $false, $true, 1, $false, "24+38=62", 2) | Out-Null
$d.Content.Find.Execute("42+11=53", $true, $true, $false, $false, $false, $true, 1, $false, "59+35=94", 2) | Out-Null
$d.Content.Find.Execute("52-28=24", $true, $true, $false, $false, $false, $true, 1, $false, "73-47=26", 2) | Out-Null
$d.Content.Find.Execute("2+63=65", $true, $true, $false, $false, $false, $true, 1, $false, "95-40=55", 2) | Out-Null
$d.Content.Find.Execute("77+21=98", $true, $true, $false, $false, $false, $true, 1, $false, "70+23=93", 2) | Out-Null
$d.Content.Find.Execute("27-17=10", $true, $true, $false, $false, $false, $true, 1, $false, "3+18=21", 2) | Out-Null
$d.Content.Find.Execute("76-45=31", $true, $true, $false, $false, $false, $true, 1, $false, "41+28=69", 2) | Out-Null
$d.Content.Find.Execute("76-38=38", $true, $true, $false, $false, $false, $true, 1, $false, "18-14=4", 2) | Out-Null
$d.Content.Find.Execute("77-41=36", $true, $true, $false, $false, $false, $true, 1, $false, "99-58=41", 2) | Out-Null
